$wb = $excel.ActiveWorkbook

# --- Sheet1: move selection from D11 to D14 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("D14").Select()

# --- Add new "Scan Path Switching" sheet right after Sheet1 ---
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "Scan Path Switching"

# --- Header row: bounding-box column names + scanpath ---
$ws2.Range("A1").Value = "min_x"
$ws2.Range("B1").Value = "min_y"
$ws2.Range("C1").Value = "min_z"
$ws2.Range("D1").Value = "max_x"
$ws2.Range("E1").Value = "max_y"
$ws2.Range("F1").Value = "max_z"

# --- Sample numeric rows ---
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Value = 3
$ws2.Range("D2").Value = 4
$ws2.Range("E2").Value = 5
$ws2.Range("F2").Value = 6

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = 3
$ws2.Range("C3").Value = 4
$ws2.Range("D3").Value = 5
$ws2.Range("E3").Value = 6
$ws2.Range("F3").Value = 7

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = 4
$ws2.Range("C4").Value = 5
$ws2.Range("D4").Value = 6
$ws2.Range("E4").Value = 7
$ws2.Range("F4").Value = 8

# --- Notes typed in column I (authoring order matters for shared-string order) ---
$ws2.Range("I3").Value = "Note: Capitalization Matters"
$ws2.Range("I4").Value = 'Note: Any unspecified areas will have `default` hatching applied to them'

# --- scanpath column header + values ---
$ws2.Range("G1").Value = "scanpath"
$ws2.Range("I2").Value = 'Eligible `scanpath` Values: `default`, `island`'
$ws2.Range("G2").Value = "island"
$ws2.Range("G3").Value = "island"
$ws2.Range("G4").Value = "frick"

# --- Column G width (widened to fit "scanpath" values) ---
$ws2.Columns.Item(7).ColumnWidth = 16.7491

# --- Final selection on the new sheet ---
$ws2.Range("G5").Select()
